$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "Wins", "Losses", "Ties" columns (AD, AE, AF),
# matching the style used by the existing header row (copy formatting from
# an existing header cell first, then set the text).
$ws.Range("A1:C1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row (2-46).
$rowCount = 45
$rng = $ws.Range("AD2:AF46")
$arr = New-Object 'object[,]' $rowCount,3
for ($i = 0; $i -lt $rowCount; $i++) {
    $arr[$i,0] = 75
    $arr[$i,1] = 87
    $arr[$i,2] = 0
}
$rng.Value = $arr
